$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06828333333333333
$ws.Range("N2").Value = 0.20485
$ws.Range("Q2").Value = 4.364068566994444
$ws.Range("R2").Value = 39.27661710295
$ws.Range("S2").Value = 0.4067926910433548
$ws.Range("T2").Value = 0.4067926910433549

# Row 3
$ws.Range("I3").Value = 0.3656254573230189
$ws.Range("J3").Value = 0.365625457323019
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06828333333333333
$ws.Range("N3").Value = 0.20485
$ws.Range("Q3").Value = 3.92242683
$ws.Range("R3").Value = 35.30184147
$ws.Range("S3").Value = 0.3656254573230189
$ws.Range("T3").Value = 0.365625457323019

# Row 4
$ws.Range("G4").Value = 35.755375
$ws.Range("H4").Value = 107.266125
$ws.Range("I4").Value = 0.2275818516336261
$ws.Range("J4").Value = 0.2275818516336262
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.06828333333333333
$ws.Range("N4").Value = 0.20485
$ws.Range("Q4").Value = 2.441496189583333
$ws.Range("R4").Value = 21.97346570625
$ws.Range("S4").Value = 0.2275818516336261
$ws.Range("T4").Value = 0.2275818516336262
